# The body text of this document is being cyclically rotated: the content
# that used to belong to one section now belongs to the section above it
# (and the very first piece of content wraps around to the last section).
# Every Find/Replace below is scoped to the Range of the specific paragraph
# (and, inside the "Avaliação" paragraph, to the sub-range that follows a
# given bold label and precedes the next one) so that look-alike text
# elsewhere in the document is never touched, and run/formatting boundaries
# (e.g. the bold "Método: " / "Critério: " / "Norma de recuperação: "
# labels, and the manual line breaks in the multi-line values) stay intact.

$d = $word.ActiveDocument

function Get-LabelRange($para, $label) {
    # Returns a Range positioned exactly on the given literal text, searched
    # for within the (still original, unmodified) paragraph.
    $r = $para.Range.Duplicate()
    $r.Find.Execute($label, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
    return $r
}

# --- "Objetivos" paragraph (Word paragraph 6) ---
# gets the old "Programa resumido" text
$r = $d.Paragraphs.Item(6).Range
$r.Find.Execute(
    "Propiciar uma integração entre os elementos de estruturação da cidade, das variáveis ambientais e da malha urbana.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Variável Ecológicano Ambiente Urbano; Enfoque Encômico e Impactos Ambientais.", 2) | Out-Null

# --- "Docente(s) Responsável(eis)" paragraph (Word paragraph 8) ---
# gets the old "Objetivos" text
$r = $d.Paragraphs.Item(8).Range
$r.Find.Execute(
    "5840942 - Marco Aurélio Kondracki de Alcântara",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Propiciar uma integração entre os elementos de estruturação da cidade, das variáveis ambientais e da malha urbana.", 2) | Out-Null

# --- "Programa resumido" paragraph (Word paragraph 10) ---
# gets the old "Programa" text
$r = $d.Paragraphs.Item(10).Range
$r.Find.Execute(
    "Variável Ecológicano Ambiente Urbano; Enfoque Encômico e Impactos Ambientais.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Elementos para estruturação ambiental da cidade. Variável ecológica no ambiente das atividades urbanas. A questão ambiental no urbanismo. A questão ambiental sob o enfoque econômico. Noções de higiene e saúde ambiental. A urbanização e os impactos ocasionados, principal enfoque da drenagem urbana.", 2) | Out-Null

# --- "Programa" paragraph (Word paragraph 12) ---
# gets the old "Método:" value from "Avaliação"
$r = $d.Paragraphs.Item(12).Range
$r.Find.Execute(
    "Elementos para estruturação ambiental da cidade. Variável ecológica no ambiente das atividades urbanas. A questão ambiental no urbanismo. A questão ambiental sob o enfoque econômico. Noções de higiene e saúde ambiental. A urbanização e os impactos ocasionados, principal enfoque da drenagem urbana.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Aula expositiva e exercícios dirigidos.", 2) | Out-Null

# --- "Avaliação" paragraph (Word paragraph 14) ---
# It holds three bold labels ("Método: ", "Critério: ", "Norma de
# recuperação: "), each immediately followed (in its own run) by a value.
# Each value is replaced with the value that used to follow the *next*
# label (and the last one gets the old Bibliografia content).
$p14 = $d.Paragraphs.Item(14)

# Método: value <- old Critério: value
$labelEnd = (Get-LabelRange $p14 "Método: ").End
$nextStart = (Get-LabelRange $p14 "Critério: ").Start
$sub = $d.Range($labelEnd, $nextStart)
$sub.Find.Execute(
    "Aula expositiva e exercícios dirigidos.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Média ponderada de exercícios e provas.", 2) | Out-Null

# Critério: value <- old Norma de recuperação: value
$labelEnd = (Get-LabelRange $p14 "Critério: ").End
$nextStart = (Get-LabelRange $p14 "Norma de recuperação: ").Start
$sub = $d.Range($labelEnd, $nextStart)
$sub.Find.Execute(
    "Média ponderada de exercícios e provas.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Prova única com nota igual ou superior a 5,0.", 2) | Out-Null

# Norma de recuperação: value <- old Bibliografia content (4 lines joined
# by manual line breaks, exactly as they appeared in the Bibliografia
# paragraph)
$labelEnd = (Get-LabelRange $p14 "Norma de recuperação: ").End
$paraEnd = $p14.Range.End
$sub = $d.Range($labelEnd, $paraEnd)
$bibText = "valle, C.R. Qualidade ambiental: o desafio de ser competitivo protegendo o meio ambiente. Pioneira. 1995.^l" + `
    "Donaire, D.. Gestão ambiental na empresa. Atlas. 2a. edição. 1999.^l" + `
    "Winter, G.. Gestão e ambiente. Modelo prático de integração empresarial. Texto Editora, Lisboa. 1992.^l" + `
    "Tucci, C.E., Porto, R.M., L.L. e Barros, M.T. org.. Drenagem Urbana. Ed. da Universidade e ABRH. 1995."
$sub.Find.Execute(
    "Prova única com nota igual ou superior a 5,0.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    $bibText, 2) | Out-Null

# --- "Bibliografia" paragraph (Word paragraph 16) ---
# gets the old "Docente(s) Responsável(eis)" content. The whole (multi-line)
# old paragraph content needs to be matched - including the manual line
# breaks ("^l") between its four "^t" runs - so every one of those runs
# gets folded into the single replacement run.
$r = $d.Paragraphs.Item(16).Range
$bibFindText = "valle, C.R. Qualidade ambiental: o desafio de ser competitivo protegendo o meio ambiente. Pioneira. 1995.^l" + `
    "Donaire, D.. Gestão ambiental na empresa. Atlas. 2a. edição. 1999.^l" + `
    "Winter, G.. Gestão e ambiente. Modelo prático de integração empresarial. Texto Editora, Lisboa. 1992.^l" + `
    "Tucci, C.E., Porto, R.M., L.L. e Barros, M.T. org.. Drenagem Urbana. Ed. da Universidade e ABRH. 1995."
$r.Find.Execute(
    $bibFindText,
    $true, $false, $false, $false, $false, $true, 1, $false,
    "5840942 - Marco Aurélio Kondracki de Alcântara", 2) | Out-Null
